# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has been
# handed off for localization (status "Ready for handoff"), with updated
# handoff file names / timestamps and an error detail message about a stale
# handback file.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab479867c834f314f33c109b2d65c9634606390a/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ccaeb25d7c16979b5d0f2fd45f68f7b946650c60/e2e/b.md."

# ---------------------------------------------------------------
# Overview sheet: row for b.md (row 3)
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-07 04:50:14"

# ---------------------------------------------------------------
# zh-cn sheet: row for b.md (row 3)
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-07 04:49:59"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Range("P1").ColumnWidth = 39.16

# ---------------------------------------------------------------
# de-de sheet: row for b.md (row 3)
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-07 04:50:14"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Range("P1").ColumnWidth = 39.16
